$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.459.11"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "3.849.19"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'715.84"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").Value = "'173.22"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "3.849.48"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "'7.34"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'36.88"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "4.490.88"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.904.20"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "71.288.54"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'7.24"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "'17.44"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "'10.77"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").Value = "'496.87"
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "'85.27"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "'10.66"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").Value = "'12.17"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'3.23"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").Value = "'2.11"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("D33").Value = "'29.55"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("D35").Value = "'9.23"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "3.808.16"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'0.104"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'6.05"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("D41").Value = "'3.37"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").Value = "'2.29"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.000320"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "'163.77"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "'48.76"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "'420.02"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "'8.65"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  -0.96%  "
